# Apply the DBC-export update to the "Autonomous_temporary" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# 1) Fill in Min/Max/Unit for the two existing RPM signals
# ---------------------------------------------------------------
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = "RPM"

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 100
$ws.Range("J15").Value = "RPM"

# ---------------------------------------------------------------
# 2) Fill in Min/Max for the IGN signal (Message: ACU_IGN)
# ---------------------------------------------------------------
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 1

# ---------------------------------------------------------------
# 3) Insert a new signal row (IGN_SWITCH) right after IGN (row 19),
#    pushing every following row down by one.
# ---------------------------------------------------------------
$ws.Rows.Item(20).Insert()

$ws.Range("A20").Value = "IGN_SWITCH"
$ws.Range("B20").Value = 8
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = "Intel"
$ws.Range("E20").Value = $false
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 1

# Match the formatting of the data row directly above it
$ws.Range("A19:K19").Copy()
$ws.Range("A20:K20").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 4) Add the Choices text for the HV signal (Message: VCU_HV),
#    which is now at row 32 after the insert above.
# ---------------------------------------------------------------
$ws.Range("K32").Value = "0=hv off, 9=hv on, 10=contactor error"

# ---------------------------------------------------------------
# 5) Append a brand-new message block: Message: RES (ID 0x191)
#    with a single signal, SIGNAL. (row 33 stays a blank separator;
#    34/35/36 are the new message header / column header / data row)
# ---------------------------------------------------------------
$ws.Range("A34").Value = "Message: RES"
$ws.Range("B34").Value = "ID: 0x191"

$ws.Range("A35").Value = "Signal Name"
$ws.Range("B35").Value = "Start Bit"
$ws.Range("C35").Value = "Length (bits)"
$ws.Range("D35").Value = "Byte Order"
$ws.Range("E35").Value = "Signed"
$ws.Range("F35").Value = "Factor"
$ws.Range("G35").Value = "Offset"
$ws.Range("H35").Value = "Min"
$ws.Range("I35").Value = "Max"
$ws.Range("J35").Value = "Unit"
$ws.Range("K35").Value = "Choices"

$ws.Range("A36").Value = "SIGNAL"
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 8
$ws.Range("D36").Value = "Intel"
$ws.Range("E36").Value = $false
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 0
$ws.Range("K36").Value = "5=GO_SIGNAL, 7=GO_SIGNAL_2, 0=EMERGENCY"

# Match formatting of the block above (Message: VCU_HV / header / HV data row,
# now at rows 30:32 after the row-20 insert). Each row shape is copied
# separately so we don't materialize stray formatted cells in columns C:K of
# the message-header row (which, like the other message rows, only uses
# columns A and B).
$ws.Range("A30:B30").Copy()
$ws.Range("A34:B34").PasteSpecial(-4122)

$ws.Range("A31:K31").Copy()
$ws.Range("A35:K35").PasteSpecial(-4122)

$ws.Range("A32:K32").Copy()
$ws.Range("A36:K36").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 6) Widen the "Choices" column (K) to fit the new longer text
# ---------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 40.2
